$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 302, shifting existing rows 302..404 down to 303..405.
$ws.Rows.Item(302).Insert()

# Populate the newly inserted row 302 with the new record's data.
$ws.Range("A302").Value = 9
$ws.Range("B302").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C302").Value = "Metropolitana"
$ws.Range("D302").Value = 45215
$ws.Range("E302").Value = 13
$ws.Range("F302").Value = 100112026
$ws.Range("G302").Value = "Haba"
$ws.Range("H302").Value = "Sin especificar"
$ws.Range("I302").Value = "Primera"
$ws.Range("J302").Value = 70
$ws.Range("K302").Value = 6000
$ws.Range("L302").Value = 6000
$ws.Range("M302").Value = 6000
$ws.Range("N302").Value = "`$/caja 20 kilos"
$ws.Range("O302").Value = "Provincia de Melipilla"
$ws.Range("P302").Value = 300
$ws.Range("Q302").Value = 20
$ws.Range("R302").Value = "Hortaliza"
